$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: new bibliography entry (article #20) ---
$ws.Range("C22").Value2 = "Development of dual inhibitors against Alzheimer’s disease using`nfragment-based QSAR and molecular docking"
$ws.Range("D22").Value2 = 2014
$ws.Range("I22").Value2 = "Goyal, M.; Dhanjal, J.K.; Goyal, S.; Tyagi, C.; Hamid, R.; Grover,A. Development of dual inhibitors against Alzheimer’s disease using fragment-based QSAR and molecular docking. BioMed Res. Int., 2014, 2014, 979606. [http://dx.doi.org/10.1155/2014/979606] [PMID: 25019089]"

# --- Reduce font size of the whole "Bibliografia" column (I2:I22) to 10pt ---
$ws.Range("I2:I22").Font.Size = 10

# I3 holds a rich-text reference (different fonts per run); the cell-level
# font change above does not touch existing run-level sizes, so resize the
# italic "Molecular diversity, 21(2), 413-426." portion explicitly too.
$ws.Range("I3").Characters(146, 36).Font.Size = 10

# --- Row 22 grew taller because of the new wrapped text (autofit) ---
$ws.Rows(22).RowHeight = 93

# --- Leave selection on the last-edited cell, like the author did ---
$ws.Range("I22").Select() | Out-Null
